# Update the acquisition timestamp (取得日時) column on the "ランサーズ" sheet
# for the existing data rows (2-7) to reflect the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-06 12:49:56"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
